# Applies the edit described by the diff on sheet "绩效表" (first sheet):
#   - Row 7 becomes a brand-new order "黎文华订单".
#   - The old row 7 content (order N2023121805 / 网络药理学...) shifts down to
#     become the new row 8.
#   - Rows 9-13 get populated with five more completed order rows that used
#     to be blank template rows.
#   - Row 29's rolled-up summary cells change to reflect 7 completed items
#     instead of 1, with the H/J lookups now erroring out (#N/A) because the
#     new business-type text doesn't match every lookup key.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# ---- Row 7: new order "黎文华订单" ----
$ws.Cells.Item(7,1).Value  = "黄礼闯"
$ws.Cells.Item(7,2).Value  = 1
$ws.Cells.Item(7,3).Value  = "黎文华订单"
$ws.Cells.Item(7,4).Value  = ""
$ws.Cells.Item(7,6).Value  = 1
$ws.Cells.Item(7,7).Value  = "菌群+对应代谢产物介导+机制研究"
$ws.Cells.Item(7,8).Value  = ""
$ws.Cells.Item(7,9).Value  = "完成"
$ws.Cells.Item(7,10).Value = ""
$ws.Cells.Item(7,11).Value = ""

# ---- Row 8: previous row 7's order (N2023121805), shifted down ----
$ws.Cells.Item(8,1).Value  = "黄礼闯"
$ws.Cells.Item(8,2).Value  = 2
$ws.Cells.Item(8,3).Value  = "N2023121805"
$ws.Cells.Item(8,4).Value  = "其他业务"
$ws.Cells.Item(8,5).Value  = "无"
$ws.Cells.Item(8,6).Value  = 1
$ws.Cells.Item(8,7).Value  = "网络药理学+Mandenol与piezo1分子对接"
$ws.Cells.Item(8,8).Value  = ""
$ws.Cells.Item(8,9).Value  = "完成"
$ws.Cells.Item(8,10).Value = "(1) 网络药理学; (2) 预测化合物靶点; (3) 获取疾病或条件相关的基因集:Genecards; (4) 疾病相关基因集:PharmGKB 数据库挖掘; (5) 疾病相关基因集:DisGeNet 数据库挖掘; (6) 调控该基因的相关转录因子 (TF) 数据获取; (7) 富集分析; (8) 全自动批量分子对接"
$ws.Cells.Item(8,11).Value = 0.1

# ---- Row 9: BI2024013001 ----
$ws.Cells.Item(9,1).Value  = "黄礼闯"
$ws.Cells.Item(9,2).Value  = 3
$ws.Cells.Item(9,3).Value  = "BI2024013001"
$ws.Cells.Item(9,4).Value  = "其他业务"
$ws.Cells.Item(9,6).Value  = 1
$ws.Cells.Item(9,7).Value  = "审核业务"
$ws.Cells.Item(9,8).Value  = ""
$ws.Cells.Item(9,9).Value  = "完成"
$ws.Cells.Item(9,10).Value = ""
$ws.Cells.Item(9,11).Value = ""

# ---- Row 10: N2024020103 ----
$ws.Cells.Item(10,1).Value  = "黄礼闯"
$ws.Cells.Item(10,2).Value  = 4
$ws.Cells.Item(10,3).Value  = "N2024020103"
$ws.Cells.Item(10,4).Value  = "其他业务"
$ws.Cells.Item(10,6).Value  = 1
$ws.Cells.Item(10,7).Value  = "筛选主动脉-下腔静脉瘘ACF模型 DEGs 并功能分析"
$ws.Cells.Item(10,8).Value  = ""
$ws.Cells.Item(10,9).Value  = "完成"
$ws.Cells.Item(10,10).Value = ""
$ws.Cells.Item(10,11).Value = ""

# ---- Row 11: N2024012602 ----
$ws.Cells.Item(11,1).Value  = "黄礼闯"
$ws.Cells.Item(11,2).Value  = 5
$ws.Cells.Item(11,3).Value  = "N2024012602"
$ws.Cells.Item(11,4).Value  = "其他业务"
$ws.Cells.Item(11,5).Value  = "1.5-2分"
$ws.Cells.Item(11,6).Value  = 1
$ws.Cells.Item(11,7).Value  = "Hydroxysafflor Yellow A 与Piezo1对接"
$ws.Cells.Item(11,8).Value  = ""
$ws.Cells.Item(11,9).Value  = "完成"
$ws.Cells.Item(11,10).Value = ""
$ws.Cells.Item(11,11).Value = ""

# ---- Row 12: N2024010303 ----
$ws.Cells.Item(12,1).Value  = "黄礼闯"
$ws.Cells.Item(12,2).Value  = 6
$ws.Cells.Item(12,3).Value  = "N2024010303"
$ws.Cells.Item(12,4).Value  = "其他业务"
$ws.Cells.Item(12,6).Value  = 1
$ws.Cells.Item(12,7).Value  = "分子对接 Celogenamide A（环状肽）蛋白 SSTR2"
$ws.Cells.Item(12,8).Value  = ""
$ws.Cells.Item(12,9).Value  = "完成"
$ws.Cells.Item(12,10).Value = "(1) 分子对接肽与蛋白"
$ws.Cells.Item(12,11).Value = ""

# ---- Row 13: 20231012 ----
# "20231012" must stay a literal TEXT order code (not become the number
# 20231012) while keeping the row's original "General" cell style (s=23).
# Plain `.Value = "20231012"` gets auto-coerced to a number by the COM
# layer, and forcing the cell to Text via .NumberFormat first would fork a
# brand-new style in styles.xml. Instead, stage the text on a cell that
# already carries a Text-formatted style (E7, which this script clears to
# "" anyway), copy it, and paste-special *values only* into C13 so the
# destination keeps its own (General) style while inheriting the text type.
$ws.Cells.Item(7,5).Value = "20231012"
$ws.Cells.Item(7,5).Copy()
$ws.Cells.Item(13,3).PasteSpecial(-4163)  # xlPasteValues
$ws.Cells.Item(7,5).Value = ""
$wb.Application.CutCopyMode = $false

$ws.Cells.Item(13,1).Value  = "黄礼闯"
$ws.Cells.Item(13,2).Value  = 7
$ws.Cells.Item(13,4).Value  = "其他业务"
$ws.Cells.Item(13,6).Value  = 1
$ws.Cells.Item(13,7).Value  = "建立风险模型和作图"
$ws.Cells.Item(13,8).Value  = ""
$ws.Cells.Item(13,9).Value  = "完成"
$ws.Cells.Item(13,10).Value = ""
$ws.Cells.Item(13,11).Value = ""

# ---- Row 29 summary recompute ----
$ws.Cells.Item(29,6).Value  = 7
$ws.Cells.Item(29,8).Value  = "#N/A"
$ws.Cells.Item(29,9).Value  = "NA+0.1+NA+NA+NA+NA+NA=NA"
$ws.Cells.Item(29,10).Value = "#N/A"
